# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers in AC1:AE1 ---
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Copy the formatting (font/border/alignment) used by the other header
# cells (e.g. A1, which carries style index 1) onto the new header cells,
# so they match the existing header row styling exactly.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# --- Data rows (2-37): season record values for every team/player row ---
$lastRow = 37
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 29).Value = 58
    $ws.Cells.Item($r, 30).Value = 56
    $ws.Cells.Item($r, 31).Value = 0
}

Write-Host "Added Wins/Losses/Ties columns (AC:AE) for rows 1-$lastRow"
